$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each changed cell, force Text number format before assigning the
# new value so Excel keeps it as a literal string (matching the original
# inlineStr cell type) instead of auto-converting numeric-looking text
# (e.g. "243.87", "0.7765") into a floating point number.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.822.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.16%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.889.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.08%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7765'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.07%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.87'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3138'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07335'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.27'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.70%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08144'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7650'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.16%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.97%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.907.54'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.77%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.02'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.04%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.208'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.69%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.828.32'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.04%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.91'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.89%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.40'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.15%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007842'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.66%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.15%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.149'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.64%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.133.79'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.95%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1568'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.410'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.36%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.90'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.77%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.49%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.036'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.89%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.75%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.541'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.31%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.473'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.75%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05570'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.25%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.070'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.38%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.246'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.30%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7540'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.53%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.50%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.636'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.95%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01931'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.55%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.776'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.23%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.145.45'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +11.39%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4441'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.61%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '73.75'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.951'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.07%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8508'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.19%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9999'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.12%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.94%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.88'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.100'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.45%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.806'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.85%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.489'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.80%  '
